$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.78044558139225
$ws.Range("C2").Value = 8.224335648617688
$ws.Range("D2").Value = 14.85986859867253
$ws.Range("E2").Value = 16.26334250504844
$ws.Range("G2").Value = 30.88242887003047
$ws.Range("H2").Value = 14.65051895999997
$ws.Range("J2").Value = 9.2715037375135
$ws.Range("N2").Value = 16.77687726198332
$ws.Range("O2").Value = 22.62262974175578
$ws.Range("B3").Value = 14.23954050940539
$ws.Range("C3").Value = 7.738148651988909
$ws.Range("D3").Value = 14.79426895644785
$ws.Range("E3").Value = 16.19779271106441
$ws.Range("G3").Value = 30.80664250112963
$ws.Range("H3").Value = 14.69090892702754
$ws.Range("J3").Value = 9.278689457990106
$ws.Range("N3").Value = 16.82205832863171
$ws.Range("O3").Value = 22.66129664898792
$ws.Range("B4").Value = 13.89850151910049
$ws.Range("C4").Value = 7.423240421401375
$ws.Range("D4").Value = 14.75730251710598
$ws.Range("E4").Value = 16.16120642878973
$ws.Range("G4").Value = 30.77183785662221
$ws.Range("H4").Value = 14.71869154685317
$ws.Range("J4").Value = 9.28453756801127
$ws.Range("N4").Value = 16.85164572057771
$ws.Range("O4").Value = 22.69116905972054
$ws.Range("B5").Value = 13.75749286081038
$ws.Range("C5").Value = 7.290860983976294
$ws.Range("D5").Value = 14.74308276002726
$ws.Range("E5").Value = 16.14722895822826
$ws.Range("G5").Value = 30.76061003903406
$ws.Range("H5").Value = 14.73076175617877
$ws.Range("J5").Value = 9.287281985278442
$ws.Range("N5").Value = 16.86416779680582
$ws.Range("O5").Value = 22.70487905436391
$ws.Range("B6").Value = 13.73396229902264
$ws.Range("C6").Value = 7.268637108248277
$ws.Range("D6").Value = 14.74077289246175
$ws.Range("E6").Value = 16.14496458246872
$ws.Range("G6").Value = 30.75892426215494
$ws.Range("H6").Value = 14.73281117037732
$ws.Range("J6").Value = 9.287759516139307
$ws.Range("N6").Value = 16.86627518436296
$ws.Range("O6").Value = 22.70724824964983
$ws.Range("B7").Value = 13.89660776244422
$ws.Range("C7").Value = 7.421471408213876
$ws.Range("D7").Value = 14.75710731129192
$ws.Range("E7").Value = 16.16101413751913
$ws.Range("G7").Value = 30.77167446327971
$ws.Range("H7").Value = 14.71885130111308
$ws.Range("J7").Value = 9.284573117289447
$ws.Range("N7").Value = 16.85181271398832
$ws.Range("O7").Value = 22.69134774222247
$ws.Range("B8").Value = 14.59591301102451
$ws.Range("C8").Value = 8.060136080452557
$ws.Range("D8").Value = 14.83656969883043
$ws.Range("E8").Value = 16.23998794448362
$ws.Range("G8").Value = 30.85386842471747
$ws.Range("H8").Value = 14.66382539271003
$ws.Range("J8").Value = 9.273683440735919
$ws.Range("N8").Value = 16.79207301071811
$ws.Range("O8").Value = 22.63468688457063
$ws.Range("B9").Value = 15.8882638102363
$ws.Range("C9").Value = 9.180245164591103
$ws.Range("D9").Value = 15.01811596065581
$ws.Range("E9").Value = 16.42337036138095
$ws.Range("G9").Value = 31.10764494961053
$ws.Range("H9").Value = 14.57965495259999
$ws.Range("J9").Value = 9.263713807293087
$ws.Range("N9").Value = 16.68953689994127
$ws.Range("O9").Value = 22.57241350085322
$ws.Range("B10").Value = 16.78003484570339
$ws.Range("C10").Value = 9.920301413421551
$ws.Range("D10").Value = 15.16635903432118
$ws.Range("E10").Value = 16.57468966914013
$ws.Range("G10").Value = 31.34963637149394
$ws.Range("H10").Value = 14.53237146666486
$ws.Range("J10").Value = 9.263313400923634
$ws.Range("N10").Value = 16.62306598165971
$ws.Range("O10").Value = 22.55666115661964
$ws.Range("B11").Value = 17.17154415788506
$ws.Range("C11").Value = 10.23868685077186
$ws.Range("D11").Value = 15.23683062964879
$ws.Range("E11").Value = 16.64693920649839
$ws.Range("G11").Value = 31.47150548070529
$ws.Range("H11").Value = 14.51403870363891
$ws.Range("J11").Value = 9.264629465901896
$ws.Range("N11").Value = 16.59474159490664
$ws.Range("O11").Value = 22.55604375362152
$ws.Range("B12").Value = 17.31764735232749
$ws.Range("C12").Value = 10.3566083007763
$ws.Range("D12").Value = 15.26393505109751
$ws.Range("E12").Value = 16.67477146258235
$ws.Range("G12").Value = 31.51931934649388
$ws.Range("H12").Value = 14.50755460188001
$ws.Range("J12").Value = 9.265342622673893
$ws.Range("N12").Value = 16.58429037976111
$ws.Range("O12").Value = 22.55675320314294
$ws.Range("B13").Value = 17.28627885911992
$ws.Range("C13").Value = 10.33132962387691
$ws.Range("D13").Value = 15.25807931838524
$ws.Range("E13").Value = 16.66875654967112
$ws.Range("G13").Value = 31.50894825735326
$ws.Range("H13").Value = 14.50893067461184
$ws.Range("J13").Value = 9.265179490122117
$ws.Range("N13").Value = 16.58652903091111
$ws.Range("O13").Value = 22.55655844204915
$ws.Range("B14").Value = 17.18360791583429
$ws.Range("C14").Value = 10.24844139207395
$ws.Range("D14").Value = 15.23905224847968
$ws.Range("E14").Value = 16.64921961069662
$ws.Range("G14").Value = 31.47540597943915
$ws.Range("H14").Value = 14.51349606421415
$ws.Range("J14").Value = 9.264683837079948
$ws.Range("N14").Value = 16.59387626653396
$ws.Range("O14").Value = 22.55608320996675
$ws.Range("B15").Value = 17.12043543269453
$ws.Range("C15").Value = 10.19732518833518
$ws.Range("D15").Value = 15.22745155157539
$ws.Range("E15").Value = 16.63731371280018
$ws.Range("G15").Value = 31.45507619787373
$ws.Range("H15").Value = 14.51635219460553
$ws.Range("J15").Value = 9.264408186457494
$ws.Range("N15").Value = 16.59841240749723
$ws.Range("O15").Value = 22.55591498948487
$ws.Range("B16").Value = 16.75415351599253
$ws.Range("C16").Value = 9.899124961113342
$ws.Range("D16").Value = 15.16181309535802
$ws.Range("E16").Value = 16.57003516818148
$ws.Range("G16").Value = 31.34190661328279
$ws.Range("H16").Value = 14.53363359399803
$ws.Range("J16").Value = 9.2632574876435
$ws.Range("N16").Value = 16.62495550500029
$ws.Range("O16").Value = 22.55683342959073
$ws.Range("B17").Value = 16.52573648104669
$ws.Range("C17").Value = 9.711495561672651
$ws.Range("D17").Value = 15.12231076929051
$ws.Range("E17").Value = 16.52962377874883
$ws.Range("G17").Value = 31.27547918133782
$ws.Range("H17").Value = 14.54504977381814
$ws.Range("J17").Value = 9.262934869445534
$ws.Range("N17").Value = 16.64172854524625
$ws.Range("O17").Value = 22.55907544155521
$ws.Range("B18").Value = 16.39303000095948
$ws.Range("C18").Value = 9.601857079409729
$ws.Range("D18").Value = 15.09987688665096
$ws.Range("E18").Value = 16.50670265540504
$ws.Range("G18").Value = 31.23838328351372
$ws.Range("H18").Value = 14.55191499582975
$ws.Range("J18").Value = 9.262890372091208
$ws.Range("N18").Value = 16.65155609997055
$ws.Range("O18").Value = 22.56098130762876
$ws.Range("B19").Value = 16.34787374664838
$ws.Range("C19").Value = 9.5644407391722
$ws.Range("D19").Value = 15.09233095756806
$ws.Range("E19").Value = 16.49899786871206
$ws.Range("G19").Value = 31.22601499080811
$ws.Range("H19").Value = 14.55429074327664
$ws.Range("J19").Value = 9.262899552068056
$ws.Range("N19").Value = 16.6549144985677
$ws.Range("O19").Value = 22.56173239389413
$ws.Range("B20").Value = 16.55019010019271
$ws.Range("C20").Value = 9.731647083310166
$ws.Range("D20").Value = 15.12648630405511
$ws.Range("E20").Value = 16.53389239620378
$ws.Range("G20").Value = 31.28243567498751
$ws.Range("H20").Value = 14.54380355326476
$ws.Range("J20").Value = 9.26295461747916
$ws.Range("N20").Value = 16.63992438615076
$ws.Range("O20").Value = 22.5587729747699
$ws.Range("B21").Value = 17.21382414230726
$ws.Range("C21").Value = 10.27285954756408
$ws.Range("D21").Value = 15.24462975390511
$ws.Range("E21").Value = 16.65494539550838
$ws.Range("G21").Value = 31.48521323927341
$ws.Range("H21").Value = 14.51214265554484
$ws.Range("J21").Value = 9.264823598695942
$ws.Range("N21").Value = 16.59171075664489
$ws.Range("O21").Value = 22.55619718876976
$ws.Range("B22").Value = 17.63495505878879
$ws.Range("C22").Value = 10.61115633830986
$ws.Range("D22").Value = 15.32427293241146
$ws.Range("E22").Value = 16.73680818464016
$ws.Range("G22").Value = 31.62742655875716
$ws.Range("H22").Value = 14.4941212197297
$ws.Range("J22").Value = 9.267296681823051
$ws.Range("N22").Value = 16.56180079066063
$ws.Range("O22").Value = 22.56001192809522
$ws.Range("B23").Value = 17.41137638136497
$ws.Range("C23").Value = 10.4320159251255
$ws.Range("D23").Value = 15.28154972608546
$ws.Range("E23").Value = 16.69287120874648
$ws.Range("G23").Value = 31.55064923066283
$ws.Range("H23").Value = 14.50349481744096
$ws.Range("J23").Value = 9.265862470709338
$ws.Range("N23").Value = 16.57761804008021
$ws.Range("O23").Value = 22.55747251787659
$ws.Range("B24").Value = 16.53913892854065
$ws.Range("C24").Value = 9.722542085507753
$ws.Range("D24").Value = 15.12459768206602
$ws.Range("E24").Value = 16.53196158130637
$ws.Range("G24").Value = 31.279287233806
$ws.Range("H24").Value = 14.54436602921576
$ws.Range("J24").Value = 9.262945250204535
$ws.Range("N24").Value = 16.6407394716762
$ws.Range("O24").Value = 22.55890779846685
$ws.Range("B25").Value = 15.54816712961215
$ws.Range("C25").Value = 8.891659224694569
$ws.Range("D25").Value = 14.96632980909228
$ws.Range("E25").Value = 16.37078558254269
$ws.Range("G25").Value = 31.0291559028961
$ws.Range("H25").Value = 14.59987494445777
$ws.Range("J25").Value = 9.263713807293087
$ws.Range("N25").Value = 16.71571605418506
$ws.Range("O25").Value = 22.58400553575599
